# Update the dSF (column F) values as part of a repull of source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -5
    "F3"  = -4
    "F4"  = -3
    "F5"  = -5
    "F7"  = -4
    "F9"  = -3
    "F11" = -3
    "F13" = -5
    "F14" = 7
    "F15" = -2
    "F16" = -1
    "F17" = -2
    "F18" = 2
    "F22" = -4
    "F24" = -1
    "F26" = 7
    "F32" = -4
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
